$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156, shifting existing rows 156..256 down to 157..257
$ws.Rows(156).EntireRow.Insert()

# Populate the new row 156 with the new record (date 2023-08-23 -> serial 45161)
$ws.Cells.Item(156, 1).Value = 11
$ws.Cells.Item(156, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(156, 3).Value = "Bíobío"
$ws.Cells.Item(156, 4).Value = 45161
$ws.Cells.Item(156, 5).Value = 8
$ws.Cells.Item(156, 6).Value = "Fruta"
$ws.Cells.Item(156, 7).Value = 100102
$ws.Cells.Item(156, 8).Value = "Cítricos"
$ws.Cells.Item(156, 9).Value = 100102004
$ws.Cells.Item(156, 10).Value = "Mandarina"
$ws.Cells.Item(156, 11).Value = "Murcott"
$ws.Cells.Item(156, 12).Value = "Primera"
$ws.Cells.Item(156, 13).Value = 270
$ws.Cells.Item(156, 14).Value = 9000
$ws.Cells.Item(156, 15).Value = 10000
$ws.Cells.Item(156, 16).Value = 9444
$ws.Cells.Item(156, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(156, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(156, 19).Value = 525
$ws.Cells.Item(156, 20).Value = 18
